$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 164.92857
$ws.Range("I33").Value = 162.23077
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 162.23077
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 66.76922999999999
$ws.Range("N33").Value = -658
$ws.Range("H53").Value = 382.81818
$ws.Range("I53").Value = 273.33334
$ws.Range("J53").Value = 514.2
$ws.Range("K53").Value = 273.33334
$ws.Range("L53").Value = 514.2
$ws.Range("M53").Value = 363.66666
$ws.Range("N53").Value = -1788.2
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H98").Value = 2600
$ws.Range("I98").Value = 2600
$ws.Range("K98").Value = 2600
$ws.Range("M98").Value = -1102
$ws.Range("H111").Value = 6516.3687
$ws.Range("J111").Value = 5994.6
$ws.Range("L111").Value = 17983.8
$ws.Range("N111").Value = -24117.8
$ws.Range("H122").Value = 2600
$ws.Range("I122").Value = 2600
$ws.Range("K122").Value = 7800
$ws.Range("M122").Value = -5350
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 686.6667
$ws.Range("I2").Value = 686.6667
$ws.Range("K2").Value = 686.6667
$ws.Range("M2").Value = -573.6667
$ws.Range("H74").Value = 1777.7
$ws.Range("I74").Value = 1959.5
$ws.Range("K74").Value = 1959.5
$ws.Range("M74").Value = -1085.5
$ws.Range("H77").Value = 1777.7
$ws.Range("I77").Value = 1959.5
$ws.Range("K77").Value = 9797.5
$ws.Range("M77").Value = -5429.5
$ws.Range("H116").Value = 686.6667
$ws.Range("I116").Value = 686.6667
$ws.Range("K116").Value = 686.6667
$ws.Range("M116").Value = 1607.3333
$ws.Range("H122").Value = 1083.3334
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1970
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 686.6667
$ws.Range("I3").Value = 686.6667
$ws.Range("K3").Value = 686.6667
$ws.Range("M3").Value = -572.6667
$ws.Range("H22").Value = 253.33333
$ws.Range("I22").Value = 253.33333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 253.33333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -80.33332999999999
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 2904.5
$ws.Range("I94").Value = 2904.5
$ws.Range("K94").Value = 2904.5
$ws.Range("M94").Value = -2453.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 1964.2
$ws.Range("I134").Value = 1964.2
$ws.Range("K134").Value = 5892.6
$ws.Range("M134").Value = -3357.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 19869
$ws.Range("J92").Value = 19869
$ws.Range("L92").Value = 19869
$ws.Range("N92").Value = -24861
$ws.Range("H107").Value = 1300.9166
$ws.Range("I107").Value = 962.25
$ws.Range("J107").Value = 1978.25
$ws.Range("K107").Value = 962.25
$ws.Range("L107").Value = 1978.25
$ws.Range("M107").Value = 957.75
$ws.Range("N107").Value = -5818.25
$ws.Range("H122").Value = 1433.7142
$ws.Range("I122").Value = 1309.5
$ws.Range("J122").Value = 1599.3334
$ws.Range("K122").Value = 3928.5
$ws.Range("L122").Value = 4798.0002
$ws.Range("M122").Value = -1478.5
$ws.Range("N122").Value = -9698.0002
$ws.Range("H132").Value = 5493.75
$ws.Range("I132").Value = 5992.6665
$ws.Range("K132").Value = 17977.9995
$ws.Range("M132").Value = -15447.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 432.66666
$ws.Range("J2").Value = 432.66666
$ws.Range("L2").Value = 2595.99996
$ws.Range("N2").Value = -2821.99996
$ws.Range("H3").Value = 8283
$ws.Range("I3").Value = 8283
$ws.Range("K3").Value = 24849
$ws.Range("M3").Value = -24737
$ws.Range("H8").Value = 108.333336
$ws.Range("I8").Value = 108.333336
$ws.Range("K8").Value = 325.000008
$ws.Range("M8").Value = -186.000008
$ws.Range("H98").Value = 3998.6667
$ws.Range("J98").Value = 4001
$ws.Range("L98").Value = 12003
$ws.Range("N98").Value = -14999
$ws.Range("H113").Value = 1498.2858
$ws.Range("J113").Value = 1438.6
$ws.Range("L113").Value = 4315.799999999999
$ws.Range("N113").Value = -8655.799999999999
$ws.Range("H132").Value = 1850
$ws.Range("I132").Value = 1850
$ws.Range("K132").Value = 16650
$ws.Range("M132").Value = -14120
$ws.Range("H133").Value = 22508.25
$ws.Range("J133").Value = 23580.857
$ws.Range("L133").Value = 70742.571
$ws.Range("N133").Value = -80862.571
$ws.Range("H140").Value = 422.25
$ws.Range("I140").Value = 422.25
$ws.Range("K140").Value = 1266.75
$ws.Range("M140").Value = 3913.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2108.6667
$ws.Range("I102").Value = 2495.5
$ws.Range("K102").Value = 2495.5
$ws.Range("M102").Value = -873.5
$ws.Range("H122").Value = 2817.111
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 6889.577
$ws.Range("I132").Value = 7339.952
$ws.Range("K132").Value = 22019.856
$ws.Range("M132").Value = -19489.856
$ws.Range("H134").Value = 32886.832
$ws.Range("J134").Value = 32886.832
$ws.Range("L134").Value = 98660.49600000001
$ws.Range("N134").Value = -103730.496
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4498.9
$ws.Range("I40").Value = 3855.5715
$ws.Range("K40").Value = 3855.5715
$ws.Range("M40").Value = -3719.5715
$ws.Range("H46").Value = 1535.7059
$ws.Range("I46").Value = 1175.7059
$ws.Range("J46").Value = 1895.7059
$ws.Range("K46").Value = 1175.7059
$ws.Range("L46").Value = 1895.7059
$ws.Range("M46").Value = -987.7058999999999
$ws.Range("N46").Value = -2271.7059
$ws.Range("H61").Value = 1859.8
$ws.Range("I61").Value = 1859.8
$ws.Range("K61").Value = 1859.8
$ws.Range("M61").Value = -1657.8
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H74").Value = 49500
$ws.Range("I74").Value = 49500
$ws.Range("K74").Value = 49500
$ws.Range("M74").Value = -48502
$ws.Range("H77").Value = 49500
$ws.Range("I77").Value = 49500
$ws.Range("K77").Value = 148500
$ws.Range("M77").Value = -143508
$ws.Range("H100").Value = 3793.9473
$ws.Range("I100").Value = 3879.2666
$ws.Range("J100").Value = 3474
$ws.Range("K100").Value = 3879.2666
$ws.Range("L100").Value = 3474
$ws.Range("M100").Value = -3338.2666
$ws.Range("N100").Value = -4556
$ws.Range("H113").Value = 1859.8
$ws.Range("I113").Value = 1859.8
$ws.Range("K113").Value = 1859.8
$ws.Range("M113").Value = 310.2
$ws.Range("H122").Value = 8891.454
$ws.Range("I122").Value = 10803.625
$ws.Range("J122").Value = 7798.7856
$ws.Range("K122").Value = 32410.875
$ws.Range("L122").Value = 23396.3568
$ws.Range("M122").Value = -29960.875
$ws.Range("N122").Value = -28296.3568
$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1600
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 4800
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = -2250
$ws.Range("N136").Value = -8400
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H126").Value = 2671.6667
$ws.Range("I126").Value = 2024.1666
$ws.Range("J126").Value = 3966.6667
$ws.Range("K126").Value = 6072.4998
$ws.Range("L126").Value = 11900.0001
$ws.Range("M126").Value = -3602.4998
$ws.Range("N126").Value = -16840.0001
$ws.Range("H136").Value = 3885.5
$ws.Range("I136").Value = 3514
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 10542
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -7992
$ws.Range("N136").Value = -20100
